$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws2 = $wb.Worksheets.Item("Test Case Steps")

# Set values first
$ws1.Range("C2:C5").Value2 = "N"

$ws1.Range("A6").Value2 = "ApplicationLinksValidationTest"
$ws1.Range("B6").Value2 = "Validate below Application links `n1. Web of Science`n2.End Note`n3.InCities`n4.ScholarOne Abstracts`n5.ScholarOne Manuscripts"
$ws1.Range("C6").Value2 = "N"
$ws1.Range("D6").Value2 = "SKIP"

$ws1.Range("A7").Value2 = "AppHeaderFooterLinkValidationTest"
$ws1.Range("B7").Value2 = "Validate Project Neon Header and Fooler links`n1.Help`n2.Cookie Policy`n3.Privacy Statement`n4.Terms of Use"
$ws1.Range("C7").Value2 = "Y"
$ws1.Range("D7").Value2 = "PASS"

# Apply formats (copy from existing cells with desired style)
$ws1.Range("A5").Copy()
$ws1.Range("A6").PasteSpecial(-4122)
$ws1.Range("D5").Copy()
$ws1.Range("D6").PasteSpecial(-4122)
$ws1.Range("A5").Copy()
$ws1.Range("C6").PasteSpecial(-4122)

$ws1.Range("A5").Copy()
$ws1.Range("A7").PasteSpecial(-4122)
$ws1.Range("D5").Copy()
$ws1.Range("D7").PasteSpecial(-4122)
$ws1.Range("C5").Copy()
$ws1.Range("C7").PasteSpecial(-4122)

$ws2.Range("A2").Copy()
$ws1.Range("B6").PasteSpecial(-4122)
$ws1.Range("B7").PasteSpecial(-4122)

# Row heights
$ws1.Rows.Item(6).RowHeight = 90
$ws1.Rows.Item(7).RowHeight = 75

# Selection
$ws1.Range("C7").Select() | Out-Null

Write-Host "done"
